$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 137.1
$ws.Range("I11").Value = 137.1
$ws.Range("K11").Value = 137.1
$ws.Range("M11").Value = 2.900000000000006

$ws.Range("H31").Value = 300
$ws.Range("J31").Value = 200
$ws.Range("L31").Value = 600
$ws.Range("N31").Value = -1060

$ws.Range("H125").Value = 2739.6667
$ws.Range("I125").Value = 532
$ws.Range("J125").Value = 3843.5
$ws.Range("K125").Value = 4788
$ws.Range("L125").Value = 34591.5
$ws.Range("M125").Value = -2328
$ws.Range("N125").Value = -39511.5

$ws.Range("H126").Value = 99500
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 99500
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 99500
$ws.Range("N126").Value = -109380

$ws.Range("H127").Value = 958.2
$ws.Range("I127").Value = 768.5
$ws.Range("J127").Value = 1717
$ws.Range("K127").Value = 2305.5
$ws.Range("L127").Value = 5151
$ws.Range("M127").Value = 2654.5
$ws.Range("N127").Value = -15071

$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0

$ws.Range("H129").Value = 1409.125
$ws.Range("I129").Value = 1354.6
$ws.Range("J129").Value = 1500
$ws.Range("K129").Value = 4063.8
$ws.Range("L129").Value = 4500
$ws.Range("M129").Value = 936.2000000000003
$ws.Range("N129").Value = -14500

$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0

$ws.Range("H131").Value = 5477.6
$ws.Range("I131").Value = 5477.6
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 16432.8
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -11392.8

$ws.Range("H132").Value = 917.0833
$ws.Range("I132").Value = 924.65
$ws.Range("J132").Value = 879.25
$ws.Range("K132").Value = 2773.95
$ws.Range("L132").Value = 2637.75
$ws.Range("M132").Value = -243.9499999999998
$ws.Range("N132").Value = -7697.75

$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0

$ws.Range("H134").Value = 124949.5
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 124949.5
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 124949.5
$ws.Range("N134").Value = -135089.5

$ws.Range("H135").Value = 889.5
$ws.Range("I135").Value = 790.5714
$ws.Range("J135").Value = 1582
$ws.Range("K135").Value = 7115.1426
$ws.Range("L135").Value = 14238
$ws.Range("M135").Value = -4580.1426
$ws.Range("N135").Value = -19308

$ws.Range("H136").Value = 125000
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 125000
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 125000
$ws.Range("N136").Value = -135200

$ws.Range("H137").Value = 3305.2415
$ws.Range("I137").Value = 1815.2354
$ws.Range("J137").Value = 5416.0835
$ws.Range("K137").Value = 5445.706200000001
$ws.Range("L137").Value = 16248.2505
$ws.Range("M137").Value = -2895.706200000001
$ws.Range("N137").Value = -21348.2505

$ws.Range("H138").Value = 5393.4116
$ws.Range("I138").Value = 1660.5834
$ws.Range("J138").Value = 14352.2
$ws.Range("K138").Value = 4981.7502
$ws.Range("L138").Value = 43056.60000000001
$ws.Range("M138").Value = 158.2497999999996
$ws.Range("N138").Value = -53336.60000000001

$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0

$ws.Range("H140").Value = 106889
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 106889
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 106889
$ws.Range("N140").Value = -117249

$ws.Range("H141").Value = 2087.2593
$ws.Range("I141").Value = 1865.625
$ws.Range("J141").Value = 3860.3333
$ws.Range("K141").Value = 5596.875
$ws.Range("L141").Value = 11580.9999
$ws.Range("M141").Value = -416.875
$ws.Range("N141").Value = -21940.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3290.4324
$ws.Range("I32").Value = 2657.0605
$ws.Range("K32").Value = 2657.0605
$ws.Range("M32").Value = -2370.0605

$ws.Range("H61").Value = 2873.5
$ws.Range("I61").Value = 2750
$ws.Range("K61").Value = 2750
$ws.Range("M61").Value = -2538

$ws.Range("H63").Value = 3830
$ws.Range("I63").Value = 500
$ws.Range("J63").Value = 5495
$ws.Range("K63").Value = 500
$ws.Range("L63").Value = 5495
$ws.Range("M63").Value = 186
$ws.Range("N63").Value = -6867

$ws.Range("H66").Value = 3830
$ws.Range("I66").Value = 500
$ws.Range("J66").Value = 5495
$ws.Range("K66").Value = 2500
$ws.Range("L66").Value = 27475
$ws.Range("M66").Value = 932
$ws.Range("N66").Value = -34339

$ws.Range("H122").Value = 920.3333
$ws.Range("I122").Value = 920.3333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2760.9999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -310.9998999999998
$ws.Range("N122").ClearContents()

$ws.Range("H136").Value = 2873.5
$ws.Range("I136").Value = 2750
$ws.Range("K136").Value = 8250
$ws.Range("M136").Value = -5700

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2902.65
$ws.Range("I31").Value = 1902.75
$ws.Range("K31").Value = 1902.75
$ws.Range("M31").Value = -1607.75

$ws.Range("H34").Value = 2902.65
$ws.Range("I34").Value = 1902.75
$ws.Range("K34").Value = 1902.75
$ws.Range("M34").Value = -1700.75

$ws.Range("H141").Value = 150000
$ws.Range("J141").Value = 150000
$ws.Range("L141").Value = 150000
$ws.Range("N141").Value = -160360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 472.16666
$ws.Range("I8").Value = 472.16666
$ws.Range("K8").Value = 1416.49998
$ws.Range("M8").Value = -1277.49998

$ws.Range("H17").Value = 99
$ws.Range("I17").Value = 99
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 297
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -128
$ws.Range("N17").ClearContents()

$ws.Range("H109").Value = 1333.3334
$ws.Range("J109").Value = 2000
$ws.Range("L109").Value = 6000
$ws.Range("N109").Value = -8080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 707.5
$ws.Range("I122").Value = 690.36365
$ws.Range("K122").Value = 2071.09095
$ws.Range("M122").Value = 378.9090500000002

$ws.Range("H132").Value = 2481
$ws.Range("I132").Value = 2166.25
$ws.Range("K132").Value = 6498.75
$ws.Range("M132").Value = -3968.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 10000
$ws.Range("J46").Value = 10000
$ws.Range("L46").Value = 10000
$ws.Range("N46").Value = -10376

$ws.Range("H100").Value = 1816.6666
$ws.Range("I100").Value = 1780
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1780
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -1239
$ws.Range("N100").Value = -3082

$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0

$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0

$ws.Range("H127").Value = 69857.5
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 69857.5
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 69857.5
$ws.Range("N127").Value = -79777.5

$ws.Range("H128").Value = 200429
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 200429
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 200429
$ws.Range("N128").Value = -210389

$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 0

$ws.Range("H130").Value = 35000
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 35000
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 35000
$ws.Range("N130").Value = -45040

$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0

$ws.Range("H132").Value = 2965.4
$ws.Range("I132").Value = 2220.7778
$ws.Range("J132").Value = 4082.3333
$ws.Range("K132").Value = 6662.3334
$ws.Range("L132").Value = 12246.9999
$ws.Range("M132").Value = -4132.3334
$ws.Range("N132").Value = -17306.9999

$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0

$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0

$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0

$ws.Range("H136").Value = 4021.9412
$ws.Range("I136").Value = 4023.3125
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 12069.9375
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -9519.9375
$ws.Range("N136").Value = -17100

$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0

$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0

$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0

$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0

$ws.Range("H141").Value = 107500
$ws.Range("I141").Value = 90000
$ws.Range("J141").Value = 125000
$ws.Range("K141").Value = 90000
$ws.Range("L141").Value = 125000
$ws.Range("M141").Value = -84820
$ws.Range("N141").Value = -135360
